# Applies the updated vm_pu.xlsx "Case_1_90" bus-voltage results (380 kV case)
# B column holds the slack-bus setpoint (1.05 -> 1.02 p.u.); C..N (skipping G/H) hold
# the recomputed per-bus voltage magnitudes for rows 2-25 (time steps 0-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"=1.02; "C2"=1.016059960140154; "D2"=1.043808470565055; "E2"=1.017584678223098; "F2"=1.048010076926771; "I2"=1.036556810286518; "J2"=1.021281926496816; "K2"=1.046581292889503; "L2"=1.020433432074942; "M2"=1.050771116642107; "N2"=1.011035789920766
    "B3"=1.02; "C3"=1.017534719921296; "D3"=1.044433685686065; "E3"=1.018853068058077; "F3"=1.048918078645702; "I3"=1.036678097909484; "J3"=1.022389206574995; "K3"=1.047018104045995; "L3"=1.021506171554783; "M3"=1.051490828883245; "N3"=1.011417396914446
    "B4"=1.02; "C4"=1.018488054533345; "D4"=1.044837602597406; "E4"=1.019673230345743; "F4"=1.049504927832158; "I4"=1.036754991276179; "J4"=1.023104433441564; "K4"=1.04729942144768; "L4"=1.022199197271131; "M4"=1.051955212944893; "N4"=1.011663483320982
    "B5"=1.02; "C5"=1.018888619627443; "D5"=1.045007255666492; "E5"=1.020017894655405; "F5"=1.049751474904736; "I5"=1.036786936937752; "J5"=1.02340481837574; "K5"=1.047417368993308; "L5"=1.022490284196846; "M5"=1.05215012510729; "N5"=1.01176673862688
    "B6"=1.02; "C6"=1.018955863752314; "D6"=1.045035732138254; "E6"=1.020075757710805; "F6"=1.049792861622224; "I6"=1.036792278448053; "J6"=1.023455237077219; "K6"=1.047437154234057; "L6"=1.022539143786387; "M6"=1.052182833268611; "N6"=1.01178406399354
    "B7"=1.02; "C7"=1.018493407747197; "D7"=1.044839870114369; "E7"=1.019677836279616; "F7"=1.049508222851284; "I7"=1.03675541963027; "J7"=1.023108448358911; "K7"=1.047300998719531; "L7"=1.02220308780948; "M7"=1.051957818607378; "N7"=1.01166486380529
    "B8"=1.02; "C8"=1.016558559472245; "D8"=1.044019896909769; "E8"=1.018013457533117; "F8"=1.048317083275701; "I8"=1.036598128672683; "J8"=1.021656400096667; "K8"=1.046729190240671; "L8"=1.020796201501807; "M8"=1.051014620010993; "N8"=1.011164930322507
    "B9"=1.02; "C9"=1.013141657744842; "D9"=1.042570130368226; "E9"=1.015076026051515; "F9"=1.046212843516564; "I9"=1.036308808434518; "J9"=1.019087869655539; "K9"=1.045711423447517; "L9"=1.018308406165053; "M9"=1.049342461692808; "N9"=1.010277496095139
    "B10"=1.02; "C10"=1.010858289757794; "D10"=1.041600377400074; "E10"=1.013114358091413; "F10"=1.044806414051939; "I10"=1.036107765727129; "J10"=1.017368613403267; "K10"=1.045026085099039; "L10"=1.016643767878419; "M10"=1.048220840625339; "N10"=1.009681422165564
    "B11"=1.02; "C11"=1.009868182853495; "D11"=1.041179700505726; "E11"=1.01226406541504; "F11"=1.044196549853722; "I11"=1.036018777389464; "J11"=1.016622454288477; "K11"=1.044727707666416; "L11"=1.015921451567432; "M11"=1.047733531050937; "N11"=1.00942223964568
    "B12"=1.02; "C12"=1.009500195492618; "D12"=1.041023327443729; "E12"=1.01194809117168; "F12"=1.043969887273559; "I12"=1.035985432421321; "J12"=1.016345035016807; "K12"=1.044616633375144; "L12"=1.015652917911042; "M12"=1.047552275088091; "N12"=1.009325803891148
    "B13"=1.02; "C13"=1.009579140008738; "D13"=1.041056875191562; "E13"=1.012015875044314; "F13"=1.044018513149981; "I13"=1.035992598185881; "J13"=1.01640455436492; "K13"=1.044640470212969; "L13"=1.015710529931169; "M13"=1.047591166324089; "N13"=1.00934649713013
    "B14"=1.02; "C14"=1.009837769389561; "D14"=1.041166777002168; "E14"=1.012237949728206; "F14"=1.044177816550036; "I14"=1.036016027014669; "J14"=1.016599528113915; "K14"=1.044718531201007; "L14"=1.015899259280054; "M14"=1.047718553431432; "N14"=1.009414271588503
    "B15"=1.02; "C15"=1.009997090409875; "D15"=1.041234475915913; "E15"=1.012374758850204; "F15"=1.044275951128727; "I15"=1.036030423765277; "J15"=1.016719622914304; "K15"=1.044766594867905; "L15"=1.016015510597623; "M15"=1.047797008005704; "N15"=1.009456007906957
    "B16"=1.02; "C16"=1.010923969832807; "D16"=1.041628280202543; "E16"=1.013170770229639; "F16"=1.044846870294389; "I16"=1.036113630815494; "J16"=1.017418097019123; "K16"=1.045045853263228; "L16"=1.016691673220074; "M16"=1.04825314710869; "N16"=1.009698600379632
    "B17"=1.02; "C17"=1.01150499794447; "D17"=1.04187509802277; "E17"=1.013669848087676; "F17"=1.045204758916409; "I17"=1.036165306245355; "J17"=1.017855769612606; "K17"=1.045220590563983; "L17"=1.017115402299305; "M17"=1.048538831233365; "N17"=1.009850482209124
    "B18"=1.02; "C18"=1.011843767917538; "D18"=1.042018988610014; "E18"=1.01396086748992; "F18"=1.045413425330618; "I18"=1.036195260897102; "J18"=1.018110892045328; "K18"=1.045322355436866; "L18"=1.017362410300202; "M18"=1.048705307763442; "N18"=1.009938968365296
    "B19"=1.02; "C19"=1.011959257260429; "D19"=1.042068039006083; "E19"=1.014060083468078; "F19"=1.045484560973397; "I19"=1.036205442980971; "J19"=1.018197854456034; "K19"=1.045357028094341; "L19"=1.017446609040667; "M19"=1.048762045158194; "N19"=1.009969122280514
    "B20"=1.02; "C20"=1.011442673035782; "D20"=1.041848624465531; "E20"=1.013616310515503; "F20"=1.045166369546563; "I20"=1.036159781277112; "J20"=1.017808828550527; "K20"=1.045201859086772; "L20"=1.017069955319491; "M20"=1.048508196400251; "N20"=1.009834197490723
    "B21"=1.02; "C21"=1.009761615555364; "D21"=1.041134416828006; "E21"=1.012172558116792; "F21"=1.044130909317547; "I21"=1.036009135833922; "J21"=1.016542120497954; "K21"=1.044695550911275; "L21"=1.015843689671553; "M21"=1.047681047943906; "N21"=1.009394318236551
    "B22"=1.02; "C22"=1.008703405084414; "D22"=1.040684701342808; "E22"=1.011264014157371; "F22"=1.043479110938348; "I22"=1.035912736969382; "J22"=1.015744168976448; "K22"=1.04437580528603; "L22"=1.015071335806316; "M22"=1.047159554440358; "N22"=1.009116799947183
    "B23"=1.02; "C23"=1.009264504941952; "D23"=1.040923166926532; "E23"=1.01174572824862; "F23"=1.043824714299738; "I23"=1.035963999258554; "J23"=1.016167324337358; "K23"=1.044545442117217; "L23"=1.015480905138712; "M23"=1.047436144151666; "N23"=1.00926400817186
    "B24"=1.02; "C24"=1.011470835393261; "D24"=1.041860586955798; "E24"=1.013640502100563; "F24"=1.045183716310263; "I24"=1.036162278349638; "J24"=1.017830039704746; "K24"=1.045210323518728; "L24"=1.017090491304824; "M24"=1.048522039451851; "N24"=1.009841556175305
    "B25"=1.02; "C25"=1.014025935260263; "D25"=1.042945503438763; "E25"=1.015835996294212; "F25"=1.046757471979914; "I25"=1.036556810286518; "J25"=1.019753091639976; "K25"=1.045975744473479; "L25"=1.018952617519486; "M25"=1.049775958593385; "N25"=1.010507697067236
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
